# Reordena la tabla de estado de cuenta (B16:G32): se agrupan los
# periodos de mora por trabajador (Num. Doc.) y se ordenan de forma
# descendente dentro de cada grupo, manteniendo el mismo conjunto de
# datos que ya existia en la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("CC", "9144662",    "DIOMEDES DE JESUS TORRES HERNANDEZ", "1903", 74453, 2233590),
    @("CC", "9144662",    "DIOMEDES DE JESUS TORRES HERNANDEZ", "1902", 89344, 2233590),
    @("CC", "9144662",    "DIOMEDES DE JESUS TORRES HERNANDEZ", "1901", 89344, 2233590),
    @("CC", "9144662",    "DIOMEDES DE JESUS TORRES HERNANDEZ", "1812", 89344, 2233590),
    @("CC", "9144662",    "DIOMEDES DE JESUS TORRES HERNANDEZ", "1811", 89344, 2233590),
    @("CC", "9023326",    "ALBEIRO JOSE COMAS MARTINEZ",        "1903", 67526, 2025780),
    @("CC", "9023326",    "ALBEIRO JOSE COMAS MARTINEZ",        "1902", 81031, 2025780),
    @("CC", "1096193978", "JOHAN ARLEY GARCIA ESPARZA",         "1903", 70904, 2127120),
    @("CC", "1096193978", "JOHAN ARLEY GARCIA ESPARZA",         "1902", 85085, 2127120),
    @("CC", "1096193978", "JOHAN ARLEY GARCIA ESPARZA",         "1901", 85085, 2127120),
    @("CC", "1096193978", "JOHAN ARLEY GARCIA ESPARZA",         "1812", 85085, 2127120),
    @("CC", "1096193978", "JOHAN ARLEY GARCIA ESPARZA",         "1811", 85085, 2127120),
    @("CC", "8867001",    "ALFONSO LUIS GALARCIO FURNIELES",    "1903", 67526, 2025780),
    @("CC", "8867001",    "ALFONSO LUIS GALARCIO FURNIELES",    "1902", 81031, 2025780),
    @("CC", "8867001",    "ALFONSO LUIS GALARCIO FURNIELES",    "1901", 81031, 2025780),
    @("CC", "8867001",    "ALFONSO LUIS GALARCIO FURNIELES",    "1812", 81031, 2025780),
    @("CC", "8867001",    "ALFONSO LUIS GALARCIO FURNIELES",    "1811", 81031, 2025780)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
